$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the value to be stored as literal text, matching the source
    # data (numeric-looking strings such as prices must stay text, not
    # get reinterpreted as Excel numbers/dates).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "98.646.42"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "3.343.00"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue $ws.Range("D5") "261.48"
$ws.Range("E5").Value = "  +0.21%  "
Set-TextValue $ws.Range("D6") "646.90"
$ws.Range("E6").Value = "  +1.62%  "
Set-TextValue $ws.Range("D7") "1.52"
$ws.Range("E7").Value = "  +7.81%  "
Set-TextValue $ws.Range("D8") "0.469"
$ws.Range("E8").Value = "  +18.47%  "
Set-TextValue $ws.Range("D9") "1.05"
$ws.Range("E9").Value = "  +18.22%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "3.335.94"
$ws.Range("E11").Value = "  -2.28%  "
Set-TextValue $ws.Range("D12") "43.72"
$ws.Range("E12").Value = "  +19.89%  "
Set-TextValue $ws.Range("D13") "0.206"
$ws.Range("E13").Value = "  +2.81%  "
$ws.Range("E14").Value = "  +10.15%  "
$ws.Range("D15").Value = "98.420.28"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "3.972.90"
$ws.Range("E16").Value = "  -1.75%  "
Set-TextValue $ws.Range("D17") "5.54"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "3.345.52"
$ws.Range("E18").Value = "  -1.25%  "
Set-TextValue $ws.Range("D19") "7.23"
$ws.Range("E19").Value = "  +15.69%  "
Set-TextValue $ws.Range("D20") "16.50"
$ws.Range("E20").Value = "  +7.28%  "
Set-TextValue $ws.Range("D21") "530.68"
$ws.Range("E21").Value = "  +7.33%  "
$ws.Range("E22").Value = "  -2.07%  "
Set-TextValue $ws.Range("D23") "10.05"
$ws.Range("E23").Value = "  +5.69%  "
Set-TextValue $ws.Range("D24") "0.0000213"
$ws.Range("E24").Value = "  -1.92%  "
Set-TextValue $ws.Range("D25") "0.418"
$ws.Range("E25").Value = "  +45.70%  "
Set-TextValue $ws.Range("D26") "100.91"
$ws.Range("E26").Value = "  +11.52%  "
Set-TextValue $ws.Range("D27") "6.00"
$ws.Range("E27").Value = "  +2.41%  "
Set-TextValue $ws.Range("D28") "12.61"
$ws.Range("E28").Value = "  +3.68%  "
$ws.Range("D29").Value = "3.522.05"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +11.51%  "
$ws.Range("E31").Value = "  +0.16%  "
Set-TextValue $ws.Range("D32") "10.82"
$ws.Range("E32").Value = "  +11.38%  "
Set-TextValue $ws.Range("D33") "0.187"
$ws.Range("E33").Value = "  -4.26%  "
Set-TextValue $ws.Range("D34") "1.01"
$ws.Range("E34").Value = "  +1.29%  "
Set-TextValue $ws.Range("D35") "28.93"
$ws.Range("E35").Value = "  +2.39%  "
Set-TextValue $ws.Range("D36") "0.515"
$ws.Range("E36").Value = "  +8.04%  "
Set-TextValue $ws.Range("D37") "7.73"
$ws.Range("E37").Value = "  +3.74%  "
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("E39").Value = "  +2.43%  "
Set-TextValue $ws.Range("D40") "522.34"
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D41") "24.71"
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D42") "1.32"
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws.Range("D43") "3.85"
$ws.Range("E43").Value = "  +0.76%  "
Set-TextValue $ws.Range("D44") "3.36"
$ws.Range("E44").Value = "  -1.35%  "
Set-TextValue $ws.Range("D45") "0.804"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("E46").Value = "  -0.03%  "
Set-TextValue $ws.Range("D47") "0.0385"
$ws.Range("E47").Value = "  +18.25%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D48") "163.85"
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D49") "2.00"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D50") "49.62"
$ws.Range("E50").Value = "  +5.90%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D51") "7.59"
$ws.Range("E51").Value = "  +14.97%  "
